# Set Exis Unit to 1 and MaxlineLoad 100%
# (ExisUnits -> column F, MaxInvest -> column I on sheet "ScenarioA", rows 8-18)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

# ExisUnits (column F) -> 0 for rows 8-18
$ws.Range("F8:F18").Value = 0

# MaxInvest (column I) -> 200 for rows 8-18
$ws.Range("I8:I18").Value = 200

# Reflect the final on-screen selection: I9:I18 selected with I9 active
$ws.Range("I9:I18").Select()
